# correct arm coding for within participant designs
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column B ("arm") values that were incorrectly numbered 2/3/4 for rows that
# belong to within-participant designs -- they should all be arm 1.
$armFixRows = 6,7,11,20,21,30,35,39,43,45,65,72,76,77,79,89,90,99,101,102,104,105,109,110,111,113
foreach ($r in $armFixRows) {
    $ws.Cells.Item($r, 2).Value = 1
}

# Clear the stray "applyFill" style that had been left on column A for these
# rows -- resetting them back to the workbook's default "Normal" style.
$styleClearRows = 75,76,77,78,79,80,88,89,90,98,99,100,101,102,103,104,105,108,109,110,111,112,113
foreach ($r in $styleClearRows) {
    $ws.Cells.Item($r, 1).Style = "Normal"
}
